$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted as row 135; the existing rows
# 135-146 (all for "Pepino ensalada" at Terminal Hortofrutícola Agro
# Chillán) shift down to 136-147. Inserting the row this way also carries
# the date-format style (s="2") from row 134 down into the new D135 cell,
# matching every other row in column D.
$ws.Range("A135:R135").Insert()

# Populate the newly inserted row 135 with the new record's data:
# Mercado ID, Mercado, Región, Fecha, Codreg, Categoría ID, Categoría,
# Variedad, Calidad, Volumen, Precio mínimo, Precio máximo,
# Precio promedio ponderado, Unidad de comercialización, Origen,
# Precio $/Kg, Kg o Unidades, Clasificación
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 44461
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112043
$ws.Cells.Item(135, 7).Value = "Pepino ensalada"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 120
$ws.Cells.Item(135, 11).Value = 16000
$ws.Cells.Item(135, 12).Value = 17000
$ws.Cells.Item(135, 13).Value = 16500
$ws.Cells.Item(135, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(135, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(135, 16).Value = 275
$ws.Cells.Item(135, 17).Value = 60
$ws.Cells.Item(135, 18).Value = "Hortaliza"
